# DSA Big O Time Complexity
# Fill in the measured running times (seconds) for each algorithm
# (m2kNaive, m2kLoop, m2kList, m2kBinS) against each input size N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N = 1000
$ws.Range("D3").Value = 0.15598100000000001
$ws.Range("E3").Value = 0.069627999999999995
$ws.Range("F3").Value = 0.025322000000000001
$ws.Range("G3").Value = 0.0028839999999999899

# N = 2000
$ws.Range("D4").Value = 0.66817700000000002
$ws.Range("E4").Value = 0.23682300000000001
$ws.Range("F4").Value = 0.097408999999999996
$ws.Range("G4").Value = 0.0059579999999999902

# N = 4000
$ws.Range("D5").Value = 2.47898
$ws.Range("E5").Value = 1.034807
$ws.Range("F5").Value = 0.38969500000000001
$ws.Range("G5").Value = 0.013141999999999999

# N = 8000
$ws.Range("D6").Value = 11.108466
$ws.Range("E6").Value = 3.8548100000000001
$ws.Range("F6").Value = 1.514867
$ws.Range("G6").Value = 0.024188000000000001

# N = 16000
$ws.Range("D7").Value = 40.771323000000002
$ws.Range("E7").Value = 16.953623
$ws.Range("F7").Value = 6.6624349999999897
$ws.Range("G7").Value = 0.049621999999999999

# N = 32000
$ws.Range("D8").Value = 160.245285
$ws.Range("E8").Value = 70.535719
$ws.Range("F8").Value = 24.306916999999999
$ws.Range("G8").Value = 0.102267

# Selection moved to I7 as part of the authoring session.
$ws.Range("I7").Select() | Out-Null
